$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4: 2024-05-05, 405 visitors, 104 subscribers, 20 favorites
$ws.Cells.Item(4, 1).Value = "05/05/2024"
$ws.Cells.Item(4, 2).Value = 405
$ws.Cells.Item(4, 3).Value = 104
$ws.Cells.Item(4, 4).Value = 20

# Row 5: 2024-05-06, 423 visitors, 107 subscribers, 21 favorites
$ws.Cells.Item(5, 1).Value = "05/06/2024"
$ws.Cells.Item(5, 2).Value = 423
$ws.Cells.Item(5, 3).Value = 107
$ws.Cells.Item(5, 4).Value = 21

# Copy styles from row 3 to rows 4 and 5 to match formatting
$ws.Range("A3:D3").Copy() | Out-Null
$ws.Range("A4:D5").PasteSpecial(-4122) | Out-Null # xlPasteFormats

$ws.Range("A5").Select() | Out-Null
